$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 129484.86
$ws.Range("J17").Value = 150865.67
$ws.Range("L17").Value = 452597.01
$ws.Range("N17").Value = -452933.01
$ws.Range("H18").Value = 24186.143
$ws.Range("I18").Value = 8667.666999999999
$ws.Range("J18").Value = 35825
$ws.Range("K18").Value = 8667.666999999999
$ws.Range("L18").Value = 35825
$ws.Range("M18").Value = -8383.666999999999
$ws.Range("N18").Value = -36393
$ws.Range("H38").Value = 364.1
$ws.Range("I38").Value = 364.1
$ws.Range("K38").Value = 1092.3
$ws.Range("M38").Value = -720.3000000000002
$ws.Range("H42").Value = 319.625
$ws.Range("I42").Value = 69.75
$ws.Range("J42").Value = 569.5
$ws.Range("K42").Value = 209.25
$ws.Range("L42").Value = 1708.5
$ws.Range("M42").Value = 20.75
$ws.Range("N42").Value = -2168.5
$ws.Range("H80").Value = 647.25
$ws.Range("I80").Value = 450
$ws.Range("J80").Value = 713
$ws.Range("K80").Value = 1350
$ws.Range("L80").Value = 2139
$ws.Range("M80").Value = -352
$ws.Range("N80").Value = -4135
$ws.Range("H83").Value = 647.25
$ws.Range("I83").Value = 450
$ws.Range("J83").Value = 713
$ws.Range("K83").Value = 4050
$ws.Range("L83").Value = 6417
$ws.Range("M83").Value = 942
$ws.Range("N83").Value = -16401
$ws.Range("H87").Value = 33354
$ws.Range("J87").Value = 33354
$ws.Range("L87").Value = 33354
$ws.Range("N87").Value = -35850
$ws.Range("H90").Value = 33354
$ws.Range("J90").Value = 33354
$ws.Range("L90").Value = 100062
$ws.Range("N90").Value = -112542
$ws.Range("H98").Value = 1012.5
$ws.Range("I98").Value = 1012.5
$ws.Range("K98").Value = 1012.5
$ws.Range("M98").Value = 485.5
$ws.Range("H103").Value = 26316274
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 27778266
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 83334798
$ws.Range("M103").Value = -614
$ws.Range("N103").Value = -83335970
$ws.Range("H112").Value = 3674.2273
$ws.Range("I112").Value = 3088.75
$ws.Range("J112").Value = 3804.3333
$ws.Range("K112").Value = 9266.25
$ws.Range("L112").Value = 11412.9999
$ws.Range("M112").Value = -8158.25
$ws.Range("N112").Value = -13628.9999
$ws.Range("H116").Value = 5357.143
$ws.Range("I116").Value = 7500
$ws.Range("K116").Value = 7500
$ws.Range("M116").Value = -4058
$ws.Range("H121").Value = 9056.200000000001
$ws.Range("J121").Value = 9056.200000000001
$ws.Range("L121").Value = 27168.6
$ws.Range("N121").Value = -30662.6
$ws.Range("H122").Value = 1012.5
$ws.Range("I122").Value = 1012.5
$ws.Range("K122").Value = 3037.5
$ws.Range("M122").Value = -587.5
$ws.Range("H123").Value = 90390
$ws.Range("J123").Value = 90390
$ws.Range("L123").Value = 90390
$ws.Range("N123").Value = -100190
$ws.Range("H125").Value = 30444.637
$ws.Range("I125").Value = 41281.125
$ws.Range("J125").Value = 1547.3334
$ws.Range("K125").Value = 371530.125
$ws.Range("L125").Value = 13926.0006
$ws.Range("M125").Value = -369070.125
$ws.Range("N125").Value = -18846.0006
$ws.Range("H132").Value = 3771.389
$ws.Range("I132").Value = 1408.4166
$ws.Range("J132").Value = 8497.333000000001
$ws.Range("K132").Value = 4225.2498
$ws.Range("L132").Value = 25491.999
$ws.Range("M132").Value = -1695.2498
$ws.Range("N132").Value = -30551.999
$ws.Range("H135").Value = 56905.777
$ws.Range("I135").Value = 1179.9
$ws.Range("K135").Value = 10619.1
$ws.Range("M135").Value = -8084.1
$ws.Range("H137").Value = 3280.4285
$ws.Range("I137").Value = 1335
$ws.Range("K137").Value = 4005
$ws.Range("M137").Value = -1455
$ws.Range("H138").Value = 1937.45
$ws.Range("I138").Value = 1070.92
$ws.Range("J138").Value = 3381.6667
$ws.Range("K138").Value = 3212.76
$ws.Range("L138").Value = 10145.0001
$ws.Range("M138").Value = 1927.24
$ws.Range("N138").Value = -20425.0001

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H5").Value = 212.08333
$ws.Range("I5").Value = 191.5
$ws.Range("J5").Value = 232.66667
$ws.Range("K5").Value = 191.5
$ws.Range("L5").Value = 232.66667
$ws.Range("M5").Value = -79.5
$ws.Range("N5").Value = -456.66667
$ws.Range("H32").Value = 31766.666
$ws.Range("I32").Value = 17180.254
$ws.Range("K32").Value = 17180.254
$ws.Range("M32").Value = -16893.254
$ws.Range("H45").Value = 462143.3
$ws.Range("I45").Value = 780033.9399999999
$ws.Range("J45").Value = 2968
$ws.Range("K45").Value = 780033.9399999999
$ws.Range("L45").Value = 2968
$ws.Range("M45").Value = -779656.9399999999
$ws.Range("N45").Value = -3722
$ws.Range("H61").Value = 1256.7084
$ws.Range("I61").Value = 987.1053000000001
$ws.Range("J61").Value = 2281.2
$ws.Range("K61").Value = 987.1053000000001
$ws.Range("L61").Value = 2281.2
$ws.Range("M61").Value = -775.1053000000001
$ws.Range("N61").Value = -2705.2
$ws.Range("H74").Value = 2231.162
$ws.Range("I74").Value = 942.16
$ws.Range("J74").Value = 4916.5835
$ws.Range("K74").Value = 942.16
$ws.Range("L74").Value = 4916.5835
$ws.Range("M74").Value = -68.15999999999997
$ws.Range("N74").Value = -6664.5835
$ws.Range("H77").Value = 2231.162
$ws.Range("I77").Value = 942.16
$ws.Range("J77").Value = 4916.5835
$ws.Range("K77").Value = 4710.8
$ws.Range("L77").Value = 24582.9175
$ws.Range("M77").Value = -342.8000000000002
$ws.Range("N77").Value = -33318.9175
$ws.Range("H102").Value = 2441.111
$ws.Range("I102").Value = 2301.182
$ws.Range("J102").Value = 2661
$ws.Range("K102").Value = 2301.182
$ws.Range("L102").Value = 2661
$ws.Range("M102").Value = -679.1819999999998
$ws.Range("N102").Value = -5905
$ws.Range("H122").Value = 2349.6667
$ws.Range("I122").Value = 1594.6
$ws.Range("J122").Value = 6125
$ws.Range("K122").Value = 4783.799999999999
$ws.Range("L122").Value = 18375
$ws.Range("M122").Value = -2333.799999999999
$ws.Range("N122").Value = -23275
$ws.Range("H132").Value = 1431.9412
$ws.Range("I132").Value = 1249.4
$ws.Range("J132").Value = 1939
$ws.Range("K132").Value = 3748.2
$ws.Range("L132").Value = 5817
$ws.Range("M132").Value = -1218.2
$ws.Range("N132").Value = -10877
$ws.Range("H136").Value = 1256.7084
$ws.Range("I136").Value = 987.1053000000001
$ws.Range("J136").Value = 2281.2
$ws.Range("K136").Value = 2961.3159
$ws.Range("L136").Value = 6843.599999999999
$ws.Range("M136").Value = -411.3159000000001
$ws.Range("N136").Value = -11943.6

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value = 212.08333
$ws.Range("I4").Value = 191.5
$ws.Range("J4").Value = 232.66667
$ws.Range("K4").Value = 191.5
$ws.Range("L4").Value = 232.66667
$ws.Range("M4").Value = -76.5
$ws.Range("N4").Value = -462.66667
$ws.Range("H12").Value = 3565.875
$ws.Range("J12").Value = 4003.8572
$ws.Range("L12").Value = 4003.8572
$ws.Range("N12").Value = -4339.8572
$ws.Range("H20").Value = 8755.521000000001
$ws.Range("I20").Value = 7605.6
$ws.Range("J20").Value = 10911.625
$ws.Range("K20").Value = 7605.6
$ws.Range("L20").Value = 10911.625
$ws.Range("M20").Value = -7358.6
$ws.Range("N20").Value = -11405.625
$ws.Range("H113").Value = 5104.2856
$ws.Range("I113").Value = 5104.2856
$ws.Range("K113").Value = 5104.2856
$ws.Range("M113").Value = -2934.2856
$ws.Range("H134").Value = 1502.3529
$ws.Range("I134").Value = 1502.3529
$ws.Range("K134").Value = 4507.0587
$ws.Range("M134").Value = -1972.0587
$ws.Range("H138").Value = 99929.875
$ws.Range("J138").Value = 99929.875
$ws.Range("L138").Value = 99929.875
$ws.Range("N138").Value = -110209.875

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 20834.666
$ws.Range("I4").Value = 62500
$ws.Range("K4").Value = 62500
$ws.Range("M4").Value = -62388
$ws.Range("H7").Value = 37037372
$ws.Range("I7").Value = 71428824
$ws.Range("J7").Value = 426
$ws.Range("K7").Value = 71428824
$ws.Range("L7").Value = 426
$ws.Range("M7").Value = -71428711
$ws.Range("N7").Value = -652
$ws.Range("H17").Value = 1408.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1408.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1408.8
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -1756.8
$ws.Range("H22").Value = 345
$ws.Range("I22").Value = 230
$ws.Range("J22").Value = 1150
$ws.Range("K22").Value = 230
$ws.Range("L22").Value = 1150
$ws.Range("M22").Value = 120
$ws.Range("N22").Value = -1850
$ws.Range("H31").Value = 2750.8462
$ws.Range("I31").Value = 2015.3334
$ws.Range("J31").Value = 3381.2856
$ws.Range("K31").Value = 2015.3334
$ws.Range("L31").Value = 3381.2856
$ws.Range("M31").Value = -1720.3334
$ws.Range("N31").Value = -3971.2856
$ws.Range("H34").Value = 2750.8462
$ws.Range("I34").Value = 2015.3334
$ws.Range("J34").Value = 3381.2856
$ws.Range("K34").Value = 2015.3334
$ws.Range("L34").Value = 3381.2856
$ws.Range("M34").Value = -1813.3334
$ws.Range("N34").Value = -3785.2856
$ws.Range("H41").Value = 14147
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 14147
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 14147
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -15003
$ws.Range("H50").Value = 11249.375
$ws.Range("J50").Value = 11249.375
$ws.Range("L50").Value = 11249.375
$ws.Range("N50").Value = -12499.375
$ws.Range("H51").Value = 10999.4
$ws.Range("J51").Value = 10999.4
$ws.Range("L51").Value = 10999.4
$ws.Range("N51").Value = -12471.4
$ws.Range("H58").Value = 958.56525
$ws.Range("I58").Value = 854.35297
$ws.Range("J58").Value = 1253.8334
$ws.Range("K58").Value = 854.35297
$ws.Range("L58").Value = 1253.8334
$ws.Range("M58").Value = -651.35297
$ws.Range("N58").Value = -1659.8334
$ws.Range("H60").Value = 10840.583
$ws.Range("J60").Value = 10908.546
$ws.Range("L60").Value = 10908.546
$ws.Range("N60").Value = -11930.546
$ws.Range("H61").Value = 10999.4
$ws.Range("J61").Value = 10999.4
$ws.Range("L61").Value = 10999.4
$ws.Range("N61").Value = -11695.4
$ws.Range("H62").Value = 4083.1667
$ws.Range("I62").Value = 4125
$ws.Range("K62").Value = 4125
$ws.Range("M62").Value = -3501
$ws.Range("H65").Value = 4083.1667
$ws.Range("I65").Value = 4125
$ws.Range("K65").Value = 20625
$ws.Range("M65").Value = -17505
$ws.Range("H68").Value = 24997.273
$ws.Range("J68").Value = 24997.273
$ws.Range("L68").Value = 24997.273
$ws.Range("N68").Value = -26495.273
$ws.Range("H71").Value = 24997.273
$ws.Range("J71").Value = 24997.273
$ws.Range("L71").Value = 74991.819
$ws.Range("N71").Value = -82479.819
$ws.Range("H74").Value = 30333.334
$ws.Range("H77").Value = 30333.334
$ws.Range("H105").Value = 2150.5
$ws.Range("I105").Value = 1978.3334
$ws.Range("K105").Value = 1978.3334
$ws.Range("M105").Value = -231.3334
$ws.Range("H132").Value = 3733.875
$ws.Range("I132").Value = 3882.4285
$ws.Range("J132").Value = 2694
$ws.Range("K132").Value = 11647.2855
$ws.Range("L132").Value = 8082
$ws.Range("M132").Value = -9117.2855
$ws.Range("N132").Value = -13142
$ws.Range("H134").Value = 2836.5293
$ws.Range("I134").Value = 3030.3333
$ws.Range("J134").Value = 2371.4
$ws.Range("K134").Value = 9090.999899999999
$ws.Range("L134").Value = 7114.200000000001
$ws.Range("M134").Value = -6555.999899999999
$ws.Range("N134").Value = -12184.2
$ws.Range("H136").Value = 958.56525
$ws.Range("I136").Value = 854.35297
$ws.Range("J136").Value = 1253.8334
$ws.Range("K136").Value = 2563.05891
$ws.Range("L136").Value = 3761.5002
$ws.Range("M136").Value = -13.0589100000002
$ws.Range("N136").Value = -8861.5002

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H11").Value = 166966.33
$ws.Range("I11").Value = 91.333336
$ws.Range("K11").Value = 274.000008
$ws.Range("M11").Value = -134.000008
$ws.Range("H55").Value = 13891710
$ws.Range("I55").Value = 934.6667
$ws.Range("J55").Value = 20837098
$ws.Range("K55").Value = 2804.0001
$ws.Range("L55").Value = 62511294
$ws.Range("M55").Value = -2627.0001
$ws.Range("N55").Value = -62511648
$ws.Range("H92").Value = 2899.125
$ws.Range("J92").Value = 2896.5
$ws.Range("L92").Value = 8689.5
$ws.Range("N92").Value = -11185.5
$ws.Range("H131").Value = 35003.223
$ws.Range("I131").Value = 1300
$ws.Range("J131").Value = 39216.125
$ws.Range("K131").Value = 3900
$ws.Range("L131").Value = 117648.375
$ws.Range("M131").Value = 1140
$ws.Range("N131").Value = -127728.375
$ws.Range("H132").Value = 2017.2703
$ws.Range("I132").Value = 918.6875
$ws.Range("J132").Value = 2854.2856
$ws.Range("K132").Value = 8268.1875
$ws.Range("L132").Value = 25688.5704
$ws.Range("M132").Value = -5738.1875
$ws.Range("N132").Value = -30748.5704

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H7").Value = 34333668
$ws.Range("J7").Value = 34333668
$ws.Range("L7").Value = 34333668
$ws.Range("N7").Value = -34333892
$ws.Range("H8").Value = 34333668
$ws.Range("J8").Value = 34333668
$ws.Range("L8").Value = 34333668
$ws.Range("N8").Value = -34333946
$ws.Range("H11").Value = 69949816
$ws.Range("I11").Value = 159458290
$ws.Range("J11").Value = 332110.88
$ws.Range("K11").Value = 159458290
$ws.Range("L11").Value = 332110.88
$ws.Range("M11").Value = -159458151
$ws.Range("N11").Value = -332388.88
$ws.Range("H12").Value = 1150
$ws.Range("J12").Value = 1100
$ws.Range("L12").Value = 1100
$ws.Range("N12").Value = -1380
$ws.Range("H18").Value = 30000
$ws.Range("I18").Value = 30000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 30000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -29707
$ws.Range("N18").ClearContents()
$ws.Range("H102").Value = 2098.4062
$ws.Range("I102").Value = 1900.4445
$ws.Range("J102").Value = 3167.4
$ws.Range("K102").Value = 1900.4445
$ws.Range("L102").Value = 3167.4
$ws.Range("M102").Value = -278.4445000000001
$ws.Range("N102").Value = -6411.4
$ws.Range("H107").Value = 22729186
$ws.Range("I107").Value = 967.8570999999999
$ws.Range("J107").Value = 33335688
$ws.Range("K107").Value = 967.8570999999999
$ws.Range("L107").Value = 33335688
$ws.Range("M107").Value = 952.1429000000001
$ws.Range("N107").Value = -33339528
$ws.Range("H122").Value = 240815
$ws.Range("I122").Value = 265795.1
$ws.Range("J122").Value = 3504
$ws.Range("K122").Value = 797385.2999999999
$ws.Range("L122").Value = 10512
$ws.Range("M122").Value = -794935.2999999999
$ws.Range("N122").Value = -15412
$ws.Range("H126").Value = 3469.6155
$ws.Range("I126").Value = 3299.1
$ws.Range("K126").Value = 9897.299999999999
$ws.Range("M126").Value = -7427.299999999999

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H31").Value = 2655.0667
$ws.Range("I31").Value = 2032.8
$ws.Range("J31").Value = 3899.6
$ws.Range("K31").Value = 2032.8
$ws.Range("L31").Value = 3899.6
$ws.Range("M31").Value = -1784.8
$ws.Range("N31").Value = -4395.6
$ws.Range("H40").Value = 5712.143
$ws.Range("I40").Value = 4022.5
$ws.Range("J40").Value = 7965
$ws.Range("K40").Value = 4022.5
$ws.Range("L40").Value = 7965
$ws.Range("M40").Value = -3886.5
$ws.Range("N40").Value = -8237
$ws.Range("H82").Value = 1623.1666
$ws.Range("I82").Value = 1325.7142
$ws.Range("J82").Value = 2039.6
$ws.Range("K82").Value = 1325.7142
$ws.Range("L82").Value = 2039.6
$ws.Range("M82").Value = -964.7141999999999
$ws.Range("N82").Value = -2761.6
$ws.Range("H85").Value = 1623.1666
$ws.Range("I85").Value = 1325.7142
$ws.Range("J85").Value = 2039.6
$ws.Range("K85").Value = 1325.7142
$ws.Range("L85").Value = 2039.6
$ws.Range("M85").Value = -77.71419999999989
$ws.Range("N85").Value = -4535.6
$ws.Range("H93").Value = 50086.145
$ws.Range("I93").Value = 2654
$ws.Range("K93").Value = 2654
$ws.Range("M93").Value = -1406
$ws.Range("H132").Value = 2965.5625
$ws.Range("I132").Value = 2437.4167
$ws.Range("J132").Value = 4550
$ws.Range("K132").Value = 7312.250100000001
$ws.Range("L132").Value = 13650
$ws.Range("M132").Value = -4782.250100000001
$ws.Range("N132").Value = -18710

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H46").Value = 51000
$ws.Range("J46").Value = 51000
$ws.Range("L46").Value = 51000
$ws.Range("N46").Value = -51462
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 6340.875
$ws.Range("I122").Value = 6073
$ws.Range("J122").Value = 7501.6665
$ws.Range("K122").Value = 18219
$ws.Range("L122").Value = 22504.9995
$ws.Range("M122").Value = -15769
$ws.Range("N122").Value = -27404.9995
$ws.Range("H132").Value = 17306.137
$ws.Range("I132").Value = 15892.262
$ws.Range("J132").Value = 46997.5
$ws.Range("K132").Value = 47676.786
$ws.Range("L132").Value = 140992.5
$ws.Range("M132").Value = -45146.786
$ws.Range("N132").Value = -146052.5
$ws.Range("H134").Value = 51000
$ws.Range("J134").Value = 51000
$ws.Range("L134").Value = 153000
$ws.Range("N134").Value = -158070
$ws.Range("H136").Value = 3355.3076
$ws.Range("I136").Value = 3556.5454
$ws.Range("J136").Value = 2248.5
$ws.Range("K136").Value = 10669.6362
$ws.Range("L136").Value = 6745.5
$ws.Range("M136").Value = -8119.636200000001
$ws.Range("N136").Value = -11845.5
